$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing labels from column B into column A and add new descriptive labels
$ws.Range("A1").Value = "First Name"
$ws.Range("A2").Value = "Middle Name"
$ws.Range("A3").Value = "Last Name"
$ws.Range("A5").Value = "User Name"
$ws.Range("A6").Value = "Password"
$ws.Range("A7").Value = "Confirm Password"
$ws.Range("A4").Value = "Employee id(INT  only )"

# Adjust column widths: A wider to fit labels, B sized for values
# (values chosen so the COM layer's internal rounding lands on the
# closest achievable stored width to the target OOXML column widths)
$ws.Columns.Item(1).ColumnWidth = 30
$ws.Columns.Item(2).ColumnWidth = 12.166666666666666
